# Fruta / hortaliza, semanal
# Insert a new weekly record as row 64 (pushing existing rows 64-100 down to
# 65-101) in the "Ciruela" subconjunto sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 64:100 down by one to make room for the new weekly entry.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record's data.
$ws.Range("A64").Value = 11
$ws.Range("B64").Value = "Vega Monumental Concepción"
$ws.Range("C64").Value = "Bíobío"
$ws.Range("D64").Value = 44960
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100103
$ws.Range("H64").Value = "Frutos de hueso (carozo)"
$ws.Range("I64").Value = 100103002
$ws.Range("J64").Value = "Ciruela"
$ws.Range("K64").Value = "Fortuna"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 100
$ws.Range("N64").Value = 10000
$ws.Range("O64").Value = 11000
$ws.Range("P64").Value = 10500
$ws.Range("Q64").Value = "$/bandeja 18 kilos granel"
$ws.Range("R64").Value = "Región de O'Higgins"
$ws.Range("S64").Value = 583
$ws.Range("T64").Value = 18
